$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds price text that sometimes looks like a plain number
# (e.g. "353.21"). Force those specific cells to Text format first so
# Excel keeps them as strings instead of auto-converting to numbers,
# matching the original inlineStr cell type.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D13",
    "D14",
    "D17",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D43",
    "D44",
    "D45",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.690.33"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "2.786.99"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "353.21"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "109.21"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").Value = "0.552"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "20.03"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "7.69"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "3.222.75"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "2.798.45"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "0.929"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "51.669.08"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "7.75"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").Value = "3.19"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").Value = "69.94"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "267.35"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "26.07"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("E28").Value = "  +12.41%  "
$ws.Range("D29").Value = "10.27"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "37.29"
$ws.Range("E30").Value = "  +7.80%  "
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").Value = "  +6.73%  "
$ws.Range("D33").Value = "51.77"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").Value = "  +8.39%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.0456"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "0.0835"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "18.54"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").Value = "120.46"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").Value = "22.05"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("D46").Value = "2.128.12"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  +5.48%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "5.44"
$ws.Range("E49").Value = "  -5.23%  "
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "0.908"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").Value = "1.34"
$ws.Range("E51").Value = "  +8.47%  "
